# ERD Diagram - slide edit
#  - Remove the (empty) title slide and the third ("Discussion Topic" solo)
#    slide, keeping only the ERD diagram slide.
#  - Update several text boxes on the remaining ERD slide:
#      * Tasks.Status bullet -> varchar enum description
#      * Employees -> Colleagues: (Name split into Fname/Lname, + Phone/Email)
#      * TasksxEmployees -> TasksxColleagues
#      * Discussion Topic -> Discussions (+ Topic bullet)
#  - Resize the two text boxes that grew new bullet lines.

$EN_DASH = [char]0x2013

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Drop the title slide and the extra ERD variant slide, leaving only
#    the main ERD diagram slide behind.
# ---------------------------------------------------------------------
$p.Slides.Item(3).Delete()   # solo "Discussion Topic" slide
$p.Slides.Item(1).Delete()   # empty title slide

$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# helpers
# ---------------------------------------------------------------------
function Get-ParaIndex($textRange, [string]$containsText) {
    $parts = $textRange.Text -split "`r"
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($parts[$i].Contains($containsText)) {
            return $i + 1
        }
    }
    throw "Paragraph not found containing: $containsText"
}

function Replace-SubText($textRange, [string]$oldText, [string]$newText) {
    $full = $textRange.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -lt 0) {
        throw "Substring not found: $oldText"
    }
    $target = $textRange.Characters($idx + 1, $oldText.Length)
    $target.Text = $newText
}

function Insert-ParaAfter($textRange, [string]$afterParaContains, [string]$newParaText) {
    $paraIdx = Get-ParaIndex $textRange $afterParaContains
    $para = $textRange.Paragraphs($paraIdx, 1)
    $para.InsertAfter("`r" + $newParaText) | Out-Null
}

# ---------------------------------------------------------------------
# 2) Tasks text box: rewrite the Status bullet as a single run.
# ---------------------------------------------------------------------
$tasks = $s.Shapes.Item("TextBox 3")
$tasksTr = $tasks.TextFrame.TextRange
Replace-SubText $tasksTr "Status – int (0 - Todo, 1 - Doing, 2 ${EN_DASH} Done)" `
    "Status ${EN_DASH} varchar( To Do, Doing, Done)"

# ---------------------------------------------------------------------
# 3) Employees text box -> Colleagues:
# ---------------------------------------------------------------------
$emp = $s.Shapes.Item("TextBox 4")
$empTr = $emp.TextFrame.TextRange

Replace-SubText $empTr "Employees" "Colleagues:"
Replace-SubText $empTr "Name" "Fname"
Insert-ParaAfter $empTr "Fname  - varchar" ("Lname" + " ${EN_DASH} varchar")
Insert-ParaAfter $empTr "Department ${EN_DASH} varchar" ("Phone ${EN_DASH} Varchar")
Insert-ParaAfter $empTr "Phone ${EN_DASH} Varchar" ("Email ${EN_DASH} Varchar")

$emp.Height = 203.5687401574803

# ---------------------------------------------------------------------
# 4) TasksxEmployees text box -> TasksxColleagues
# ---------------------------------------------------------------------
$taskEmp = $s.Shapes.Item("TextBox 5")
$taskEmpTr = $taskEmp.TextFrame.TextRange
Replace-SubText $taskEmpTr "TasksxEmployees" "TasksxColleagues"

# ---------------------------------------------------------------------
# 5) Discussion Topic text box -> Discussions (+ Topic bullet)
# ---------------------------------------------------------------------
$disc = $s.Shapes.Item("TextBox 25")
$discTr = $disc.TextFrame.TextRange
Replace-SubText $discTr "Discussion Topic" "Discussions"
Insert-ParaAfter $discTr "Content ${EN_DASH} varchar" ("Topic ${EN_DASH} Varchar")

$disc.Height = 181.75779527559055
